# Apply weekly re-shuffle of Fecha (D) / Volumen (M) values, and swap the
# Unidad de comercialización / Precio $/Kg / Kg por unidad (Q/S/T) values
# between rows 4 and 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Fecha) - new serial date values
$ws.Range("D2").Value = 44322
$ws.Range("D3").Value = 44313
$ws.Range("D4").Value = 44327
$ws.Range("D5").Value = 44323
$ws.Range("D6").Value = 44316
$ws.Range("D7").Value = 44302
$ws.Range("D8").Value = 44309
$ws.Range("D9").Value = 44330
$ws.Range("D10").Value = 44306

# Column M (Volumen) - new values
$ws.Range("M2").Value = 60
$ws.Range("M4").Value = 60
$ws.Range("M5").Value = 80
$ws.Range("M6").Value = 120
$ws.Range("M7").Value = 80
$ws.Range("M9").Value = 60
$ws.Range("M10").Value = 80

# Rows 4 and 8 swap their Unidad de comercialización (Q), Precio $/Kg (S)
# and Kg / unidad (T) values.
$ws.Range("Q4").Value = "$/caja 10 kilos empedrada"
$ws.Range("S4").Value = 11500
$ws.Range("T4").Value = 1

$ws.Range("Q8").Value = "$/caja 14 kilos granel"
$ws.Range("S8").Value = 821
$ws.Range("T8").Value = 14
